# Fruta / hortaliza, semanal
# Insert 5 new weekly price rows into the "Vega Modelo de Temuco - Nectarín"
# sheet, right before the existing row 702. This pushes the former rows
# 702-708 down to 707-713 (unchanged), and the 5 freshly inserted rows
# (702-706) get populated with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows above row 702 - Excel shifts rows 702:708 down to 707:713,
# carrying their formatting (e.g. the date style on column D) with them.
$ws.Rows.Item(702).Resize(5).Insert()

# Columns A, B, C, E, F, G, H, I, J are constant for every row in this block
# (same mercado / rubro / categoría / producto), so copy them down into the
# newly-inserted rows from the row right above (701), which still holds them.
for ($r = 702; $r -le 706; $r++) {
    $ws.Cells.Item($r, 1).Value  = $ws.Cells.Item(701, 1).Value()
    $ws.Cells.Item($r, 2).Value  = $ws.Cells.Item(701, 2).Value()
    $ws.Cells.Item($r, 3).Value  = $ws.Cells.Item(701, 3).Value()
    $ws.Cells.Item($r, 5).Value  = $ws.Cells.Item(701, 5).Value()
    $ws.Cells.Item($r, 6).Value  = $ws.Cells.Item(701, 6).Value()
    $ws.Cells.Item($r, 7).Value  = $ws.Cells.Item(701, 7).Value()
    $ws.Cells.Item($r, 8).Value  = $ws.Cells.Item(701, 8).Value()
    $ws.Cells.Item($r, 9).Value  = $ws.Cells.Item(701, 9).Value()
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item(701, 10).Value()
}

# New row 702: Nectarín Artic Star, Primera, $/bandeja 18 kilos granel
$ws.Cells.Item(702, 4).Value  = 44939
$ws.Cells.Item(702, 11).Value = "Artic Star"
$ws.Cells.Item(702, 12).Value = "Primera"
$ws.Cells.Item(702, 13).Value = 300
$ws.Cells.Item(702, 14).Value = 18000
$ws.Cells.Item(702, 15).Value = 18000
$ws.Cells.Item(702, 16).Value = 18000
$ws.Cells.Item(702, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(702, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(702, 19).Value = 1000
$ws.Cells.Item(702, 20).Value = 18

# New row 703: Nectarín Artic Star, Primera, $/bins (420 kilos)
$ws.Cells.Item(703, 4).Value  = 44939
$ws.Cells.Item(703, 11).Value = "Artic Star"
$ws.Cells.Item(703, 12).Value = "Primera"
$ws.Cells.Item(703, 13).Value = 3
$ws.Cells.Item(703, 14).Value = 440000
$ws.Cells.Item(703, 15).Value = 440000
$ws.Cells.Item(703, 16).Value = 440000
$ws.Cells.Item(703, 17).Value = "$/bins (420 kilos)"
$ws.Cells.Item(703, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(703, 19).Value = 1048
$ws.Cells.Item(703, 20).Value = 420

# New row 704: Nectarín Early John, Primera, $/bandeja 18 kilos granel
$ws.Cells.Item(704, 4).Value  = 44939
$ws.Cells.Item(704, 11).Value = "Early John"
$ws.Cells.Item(704, 12).Value = "Primera"
$ws.Cells.Item(704, 13).Value = 200
$ws.Cells.Item(704, 14).Value = 18000
$ws.Cells.Item(704, 15).Value = 18000
$ws.Cells.Item(704, 16).Value = 18000
$ws.Cells.Item(704, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(704, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(704, 19).Value = 1000
$ws.Cells.Item(704, 20).Value = 18

# New row 705: Nectarín Ruby Diamond, Primera, $/bandeja 18 kilos granel
$ws.Cells.Item(705, 4).Value  = 44939
$ws.Cells.Item(705, 11).Value = "Ruby Diamond"
$ws.Cells.Item(705, 12).Value = "Primera"
$ws.Cells.Item(705, 13).Value = 250
$ws.Cells.Item(705, 14).Value = 18000
$ws.Cells.Item(705, 15).Value = 18000
$ws.Cells.Item(705, 16).Value = 18000
$ws.Cells.Item(705, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(705, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(705, 19).Value = 1000
$ws.Cells.Item(705, 20).Value = 18

# New row 706: Nectarín Ruby Diamond, Primera, $/bins (420 kilos)
$ws.Cells.Item(706, 4).Value  = 44939
$ws.Cells.Item(706, 11).Value = "Ruby Diamond"
$ws.Cells.Item(706, 12).Value = "Primera"
$ws.Cells.Item(706, 13).Value = 3
$ws.Cells.Item(706, 14).Value = 420000
$ws.Cells.Item(706, 15).Value = 420000
$ws.Cells.Item(706, 16).Value = 420000
$ws.Cells.Item(706, 17).Value = "$/bins (420 kilos)"
$ws.Cells.Item(706, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(706, 19).Value = 1000
$ws.Cells.Item(706, 20).Value = 420
